$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A72 A73 調整加入 WT & Default
$ws.Range("A19").Value = "WT"
$ws.Range("B19").Value = 0.52716666666666701
$ws.Range("A20").Value = "Default"
$ws.Range("B20").Value = 100
